# Commit: "Updated last two slides."
#
# The author removed the slide titled "Deploying to Azure with the Command
# Line" from the deck entirely. The two slides that used to follow it -
# "Using Azure Storage in a WebJob" and "Deploying and Debugging" - along
# with the trailing (title-less) slide are untouched in content; they
# simply move up one position once the slide ahead of them is gone, which
# is exactly what deleting that one slide from the slide list produces.

$p = $ppt.ActivePresentation

$targetTitle = "Deploying to Azure with the Command Line"
$targetIndex = -1

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    if ($s.Shapes.Count -gt 0 -and $s.Shapes.Item(1).HasTextFrame) {
        $text = $s.Shapes.Item(1).TextFrame.TextRange.Text
        if ($text -eq $targetTitle) {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -eq -1) {
    # Fall back to the known position in the original deck (8th slide).
    $targetIndex = 8
}

$p.Slides.Item($targetIndex).Delete()
